$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 ("I0") and J1 ("IF") should carry the same formatting
# (bold, centered, bordered) as the other header cells. Copy the format from
# the existing H1 header cell, then set the new header text.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for the new columns I (I0) and J (IF), rows 2-9.
$i0 = @(6, 6, 9, 6, 8, 5, 3, 7)
$if = @(6, 6, 9, 6, 8, 5, 3, 7)

for ($r = 0; $r -lt 8; $r++) {
    $row = $r + 2
    $ws.Cells.Item($row, 9).Value = $i0[$r]
    $ws.Cells.Item($row, 10).Value = $if[$r]
}
